$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows: "penghasilan teratur" list ---------------------------
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = 20000
$ws.Range("A3").Value = "00124"
$ws.Range("B3").Value = 10000

# --- Header row styling (A1:B1) -> accent-colored header with border -----
$a1 = $ws.Range("A1")
$a1.Interior.Color = 12874308
$a1.Font.Color = 16777215
$a1.Borders.LineStyle = 1

$b1 = $ws.Range("B1")
$b1.Interior.Color = 12874308
$b1.Font.Color = 16777215
$b1.Borders.LineStyle = 1

# --- Selection moves to C3 -------------------------------------------------
$ws.Range("C3").Select()
